# Update "Generate Report for Handback" timestamps for the
# 94dd3d62-dac4-4aad-9176-c5e8d6b1c424 row across the Overview, zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" column (G) for the
# 94dd3d62-... row.
$wsOverview.Range("G4").Value = "2016-09-03 22:48:39"

# zh-cn sheet - "Correspond Handoff Datetime" (H) / "Correspond Handback
# DateTime" (K) for the 94dd3d62-... row.
$wsZhCn.Range("H4").Value = "2016-09-03 22:48:35"
$wsZhCn.Range("K4").Value = "2016-09-03 22:48:52"

# de-de sheet - same two columns for the 94dd3d62-... row.
$wsDeDe.Range("H4").Value = "2016-09-03 22:48:39"
$wsDeDe.Range("K4").Value = "2016-09-03 22:49:00"

Write-Host "Updated handback status timestamps"
